# Edit: remove the underline formatting that split each of the three
# "X leads to Y" bullets on slide 4 into three separate runs, collapsing
# each paragraph back down into a single, plainly-formatted run (the
# formatting that survives is that of the middle " leads to " run, which
# never had the underline applied).
#
# Also best-effort touch the (normally GUI-only) slide-guide collection,
# since the canonical edit shows PowerPoint persisted an empty
# <p:extLst><p:ext ...><p15:sldGuideLst/></p:ext></p:extLst> block on
# p:presentation after the save that produced this commit.

$p = $ppt.ActivePresentation

# --- Best effort: touch presentation Guides (matches the empty
# p15:sldGuideLst extLst block PowerPoint writes out after the Guides UI
# has been touched). Harmless no-op if unsupported by the host.
try {
    $null = $p.Guides.Count
} catch {
}

$slide = $p.Slides.Item(4)
$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

function Merge-LeadsToParagraph($paraIndex, $fullText) {
    $para = $tr.Paragraphs($paraIndex, 1)

    # Run 1 = the first (underlined) phrase, Run 2 = " leads to ",
    # Run 3 = the second (underlined) phrase. Emptying a run's Text
    # deletes that run outright, so after clearing run 1 the old run 3
    # becomes run 2, etc.
    $firstRun = $para.Runs(1, 1)
    $firstRun.Text = ""

    $lastRun = $para.Runs(2, 1)
    $lastRun.Text = ""

    # Only the plain " leads to " run is left; give it the merged text.
    $remainingRun = $para.Runs(1, 1)
    $remainingRun.Text = $fullText
}

Merge-LeadsToParagraph 1 "Missing information leads to lack of content variety"
Merge-LeadsToParagraph 3 "Lack of content variety leads to bias"
Merge-LeadsToParagraph 5 "Bias leads to discrimination"
